$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5523.6313
$ws.Range("J62").Value = 5646
$ws.Range("L62").Value = 5646
$ws.Range("N62").Value = -6894
$ws.Range("H65").Value = 5523.6313
$ws.Range("J65").Value = 5646
$ws.Range("L65").Value = 28230
$ws.Range("N65").Value = -34470
$ws.Range("H94").Value = 6450
$ws.Range("I94").Value = 2900
$ws.Range("K94").Value = 2900
$ws.Range("M94").Value = -2449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 6247.5977
$ws.Range("I32").Value = 3163.0547
$ws.Range("J32").Value = 31266.666
$ws.Range("K32").Value = 3163.0547
$ws.Range("L32").Value = 31266.666
$ws.Range("M32").Value = -2876.0547
$ws.Range("N32").Value = -31840.666
$ws.Range("H52").Value = 17999.334
$ws.Range("J52").Value = 17999.334
$ws.Range("L52").Value = 17999.334
$ws.Range("N52").Value = -18635.334
$ws.Range("H88").Value = 2103.8572
$ws.Range("I88").Value = 2133.7144
$ws.Range("J88").Value = 2074
$ws.Range("K88").Value = 2133.7144
$ws.Range("L88").Value = 2074
$ws.Range("M88").Value = -1727.7144
$ws.Range("N88").Value = -2886
$ws.Range("H91").Value = 2103.8572
$ws.Range("I91").Value = 2133.7144
$ws.Range("J91").Value = 2074
$ws.Range("K91").Value = 2133.7144
$ws.Range("L91").Value = 2074
$ws.Range("M91").Value = -729.7143999999998
$ws.Range("N91").Value = -4882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H112").Value = 29411.5
$ws.Range("J112").Value = 29411.5
$ws.Range("L112").Value = 29411.5
$ws.Range("N112").Value = -32365.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.5
$ws.Range("I7").Value = 23.75
$ws.Range("J7").Value = 91.25
$ws.Range("K7").Value = 23.75
$ws.Range("L7").Value = 91.25
$ws.Range("M7").Value = 89.25
$ws.Range("N7").Value = -317.25
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650
$ws.Range("H62").Value = 2519.2727
$ws.Range("I62").Value = 2300
$ws.Range("K62").Value = 2300
$ws.Range("M62").Value = -1676
$ws.Range("H65").Value = 2519.2727
$ws.Range("I65").Value = 2300
$ws.Range("K65").Value = 11500
$ws.Range("M65").Value = -8380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1077.5555
$ws.Range("I80").Value = 866.3333
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 2598.9999
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -1662.9999
$ws.Range("N80").Value = -6372
$ws.Range("H83").Value = 1077.5555
$ws.Range("I83").Value = 866.3333
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 7796.9997
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -3116.9997
$ws.Range("N83").Value = -22860
$ws.Range("H96").Value = 4837.5
$ws.Range("J96").Value = 4837.5
$ws.Range("L96").Value = 14512.5
$ws.Range("N96").Value = -18630.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 97.09524
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 158.09091
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 158.09091
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -384.09091
$ws.Range("H19").Value = 2995
$ws.Range("I19").Value = 990
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 990
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -702
$ws.Range("N19").Value = -5576
$ws.Range("H80").Value = 2431.1365
$ws.Range("I80").Value = 2405
$ws.Range("J80").Value = 2500.8333
$ws.Range("K80").Value = 2405
$ws.Range("L80").Value = 2500.8333
$ws.Range("M80").Value = -1407
$ws.Range("N80").Value = -4496.8333
$ws.Range("H83").Value = 2431.1365
$ws.Range("I83").Value = 2405
$ws.Range("J83").Value = 2500.8333
$ws.Range("K83").Value = 12025
$ws.Range("L83").Value = 12504.1665
$ws.Range("M83").Value = -7033
$ws.Range("N83").Value = -22488.1665
$ws.Range("H107").Value = 978.35297
$ws.Range("I107").Value = 841.0909
$ws.Range("J107").Value = 1230
$ws.Range("K107").Value = 841.0909
$ws.Range("L107").Value = 1230
$ws.Range("M107").Value = 1078.9091
$ws.Range("N107").Value = -5070
$ws.Range("H132").Value = 3127721.8
$ws.Range("I132").Value = 3468.1538
$ws.Range("K132").Value = 10404.4614
$ws.Range("M132").Value = -7874.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8356.235000000001
$ws.Range("I68").Value = 26575
$ws.Range("J68").Value = 2750.4614
$ws.Range("K68").Value = 26575
$ws.Range("L68").Value = 2750.4614
$ws.Range("M68").Value = -25826
$ws.Range("N68").Value = -4248.4614
$ws.Range("H71").Value = 8356.235000000001
$ws.Range("I71").Value = 26575
$ws.Range("J71").Value = 2750.4614
$ws.Range("K71").Value = 132875
$ws.Range("L71").Value = 13752.307
$ws.Range("M71").Value = -129131
$ws.Range("N71").Value = -21240.307
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 70357.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 70357.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 70357.25
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -80557.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9533.333000000001
$ws.Range("J54").Value = 9533.333000000001
$ws.Range("L54").Value = 9533.333000000001
$ws.Range("N54").Value = -10573.333
$ws.Range("H62").Value = 3250
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3250
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4498
$ws.Range("H65").Value = 3250
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 16250
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -22490
$ws.Range("H81").Value = 1333.3334
$ws.Range("I81").Value = 1250
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 2500
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1439
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1333.3334
$ws.Range("I84").Value = 1250
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 12500
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -7196
$ws.Range("N84").Value = -25608
